$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 24 ("07/2025" / day 23), shifting the existing
# rows 24-114 down to 25-115. We shift the data manually (bottom-up copy)
# instead of relying on Rows.Insert so every source/destination pair is
# explicit and unambiguous.
for ($r = 114; $r -ge 24; $r--) {
    $dst = $r + 1
    $ws.Cells.Item($dst, 1).Value2 = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($dst, 2).Value2 = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($dst, 3).Value2 = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($dst, 4).Value2 = $ws.Cells.Item($r, 4).Value2
    $ws.Cells.Item($dst, 5).Value2 = $ws.Cells.Item($r, 5).Value2
}

# Fill in the newly freed row 24 with the new day's data (07/2025, day 23).
$ws.Cells.Item(24, 1).Value2 = 23
$ws.Cells.Item(24, 2).Value2 = 14310.32
$ws.Cells.Item(24, 3).Value2 = 7
$ws.Cells.Item(24, 4).Value2 = 2025
$ws.Cells.Item(24, 5).Value2 = "07/2025"
